# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Acelga"
# at row 97 (pushing the existing rows 97:133 down to 98:134), matching the
# near-duplicate entry already present at (old) row 111 but two days later
# (44553 instead of 44551) - "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97:133 down to 98:134, leaving a blank row 97 to fill in.
$ws.Rows("97:97").Insert()

$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 44553
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 100112009
$ws.Range("G97").Value = "Acelga"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 50
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = 10000
$ws.Range("N97").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O97").Value = "Región de La Araucanía"
$ws.Range("P97").Value = 833
$ws.Range("Q97").Value = 12
$ws.Range("R97").Value = "Hortaliza"
